# NIT-8001511751.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement table previously held a single worker / single period
# row. The edit expands it to 3 workers x 2 periods (6 data rows), updates the
# summary counters and the total "VALOR MORA", and (as a consequence of the
# extra rows) pushes the signature block further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for 5 more data rows below the existing one (row 16).
#    Inserting blank rows 17:21 shifts the blank spacer rows and the
#    signature block down by 5 (old 21/22 -> new 26/27), exactly like the
#    target layout.
# ---------------------------------------------------------------------------
$ws.Rows("17:21").Insert()

# Copy row 16's formatting (borders, fonts, number formats) down onto the
# 5 freshly-inserted rows so the new records look like the original one.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the 6 data rows (3 workers x 2 periods).
# ---------------------------------------------------------------------------
$rows = @(
  @{ r=16; tipo="CC"; doc="73189921";   nombre="TOMAS ELIAS VILORIA GUTIERREZ";    periodo="2507"; mora=28470; salario=1241560 },
  @{ r=17; tipo="CC"; doc="73189921";   nombre="TOMAS ELIAS VILORIA GUTIERREZ";    periodo="2506"; mora=28470; salario=1241560 },
  @{ r=18; tipo="CC"; doc="1047482848"; nombre="SINDY LIBETH AYAZO VILORIA";       periodo="2507"; mora=56940; salario=1423500 },
  @{ r=19; tipo="CC"; doc="1047482848"; nombre="SINDY LIBETH AYAZO VILORIA";       periodo="2506"; mora=56940; salario=1423500 },
  @{ r=20; tipo="CC"; doc="1006582607"; nombre="YEIDIS PAOLA VILLALBA BERMUDEZ";   periodo="2507"; mora=56940; salario=1423500 },
  @{ r=21; tipo="CC"; doc="1006582607"; nombre="YEIDIS PAOLA VILLALBA BERMUDEZ";   periodo="2506"; mora=22776; salario=1423500 }
)

foreach ($row in $rows) {
  $n = $row.r
  $ws.Range("B$n").Value = $row.tipo
  $ws.Range("C$n").Value = $row.doc
  $ws.Range("D$n").Value = $row.nombre
  $ws.Range("E$n").Value = $row.periodo
  $ws.Range("F$n").Value = $row.mora
  $ws.Range("G$n").Value = $row.salario
}

# ---------------------------------------------------------------------------
# 3. Update the summary block above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 250536   # VALOR MORA (sum of the 6 rows above)
$ws.Range("C13").Value = 3        # Cant. Trabajadores
$ws.Range("F13").Value = 2        # Cant. Periodos

Write-Output "Edit applied"
